# Update automatico via Actualizar 06-12-2020 05-53-22
#
# Appends the new daily record row (11/069) to the "Condicion_Pacientes"
# table on Hoja1, letting the ListObject grow so the table range /
# autofilter / dimension all expand from F90 -> F91 automatically, then
# copies the formatting from the previous row so the new cells keep the
# same styles (date-style cell for column A, centered style for B:F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$lo = $ws.ListObjects.Item("Condicion_Pacientes")

# Grow the table by one row (this updates the table ref, autofilter ref
# and worksheet dimension for us).
$newListRow = $lo.ListRows.Add()

$lastRow = $lo.Range.Rows.Count + $lo.Range.Row - 1
$newRow = $lastRow
$prevRow = $newRow - 1

# Fill in the new record.
$ws.Cells.Item($newRow, 1).Value = "11/069"
$ws.Cells.Item($newRow, 2).Value = 828
$ws.Cells.Item($newRow, 3).Value = 309
$ws.Cells.Item($newRow, 4).Value = 390
$ws.Cells.Item($newRow, 5).Value = 264
$ws.Cells.Item($newRow, 6).Value = 53

# Match the formatting of the row above (date-style for col A, centered
# for the rest) by copying its formats onto the new row.
$ws.Range($ws.Cells.Item($prevRow, 1), $ws.Cells.Item($prevRow, 6)).Copy() | Out-Null
$ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 6)).PasteSpecial(-4122) | Out-Null

# Re-apply the values in case PasteSpecial(formats) touched them.
$ws.Cells.Item($newRow, 1).Value = "11/069"
$ws.Cells.Item($newRow, 2).Value = 828
$ws.Cells.Item($newRow, 3).Value = 309
$ws.Cells.Item($newRow, 4).Value = 390
$ws.Cells.Item($newRow, 5).Value = 264
$ws.Cells.Item($newRow, 6).Value = 53

# Move the view/selection to the newly added row, mirroring what the
# author's session looked like after adding the record.
$ws.Cells.Item($newRow, 4).Select() | Out-Null

$excel.CutCopyMode = $false
